$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (station_id, model_number, available_bays, last_maintenance_date, location_id)
# is being removed; all data below shifts up by one row, matching the
# behaviour of selecting the row-1 header in Excel's UI and deleting it.
$ws.Rows(1).Select()
$ws.Rows(1).Delete()
